$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1711.9
$ws.Range("I38").Value = 24
$ws.Range("J38").Value = 3399.8
$ws.Range("K38").Value = 72
$ws.Range("L38").Value = 10199.4
$ws.Range("M38").Value = 300
$ws.Range("N38").Value = -10943.4
$ws.Range("H39").Value = 107.8
$ws.Range("I39").Value = 118
$ws.Range("J39").Value = 92.5
$ws.Range("K39").Value = 354
$ws.Range("L39").Value = 277.5
$ws.Range("M39").Value = -58
$ws.Range("N39").Value = -869.5
$ws.Range("H76").Value = 4176.923
$ws.Range("I76").Value = 4600
$ws.Range("J76").Value = 3225
$ws.Range("K76").Value = 4600
$ws.Range("L76").Value = 3225
$ws.Range("M76").Value = -4285
$ws.Range("N76").Value = -3855
$ws.Range("H79").Value = 4176.923
$ws.Range("I79").Value = 4600
$ws.Range("J79").Value = 3225
$ws.Range("K79").Value = 4600
$ws.Range("L79").Value = 3225
$ws.Range("M79").Value = -3508
$ws.Range("N79").Value = -5409
$ws.Range("H112").Value = 1993.0667
$ws.Range("J112").Value = 2317.2942
$ws.Range("L112").Value = 6951.882599999999
$ws.Range("N112").Value = -9167.882599999999
$ws.Range("H129").Value = 588.9729599999999
$ws.Range("J129").Value = 855.94446
$ws.Range("L129").Value = 2567.83338
$ws.Range("N129").Value = -12567.83338
$ws.Range("H137").Value = 1040.9324
$ws.Range("I137").Value = 749.0213
$ws.Range("J137").Value = 1549.0741
$ws.Range("K137").Value = 2247.0639
$ws.Range("L137").Value = 4647.2223
$ws.Range("M137").Value = 302.9360999999999
$ws.Range("N137").Value = -9747.222300000001
$ws.Range("H138").Value = 1412.63
$ws.Range("I138").Value = 815.4194
$ws.Range("J138").Value = 1680.942
$ws.Range("K138").Value = 2446.2582
$ws.Range("L138").Value = 5042.826
$ws.Range("M138").Value = 2693.7418
$ws.Range("N138").Value = -15322.826

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 1900
$ws.Range("I21").Value = 1900
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 1900
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -1526
$ws.Range("N21").Value = $null
$ws.Range("H28").Value = 6294.636
$ws.Range("I28").Value = 6294.636
$ws.Range("K28").Value = 6294.636
$ws.Range("M28").Value = -6102.636
$ws.Range("H61").Value = 30304240
$ws.Range("I61").Value = 40000876
$ws.Range("K61").Value = 40000876
$ws.Range("M61").Value = -40000664
$ws.Range("H74").Value = 1073.5278
$ws.Range("I74").Value = 871.55884
$ws.Range("J74").Value = 4507
$ws.Range("K74").Value = 871.55884
$ws.Range("L74").Value = 4507
$ws.Range("M74").Value = 2.441159999999968
$ws.Range("N74").Value = -6255
$ws.Range("H77").Value = 1073.5278
$ws.Range("I77").Value = 871.55884
$ws.Range("J77").Value = 4507
$ws.Range("K77").Value = 4357.7942
$ws.Range("L77").Value = 22535
$ws.Range("M77").Value = 10.20579999999973
$ws.Range("N77").Value = -31271
$ws.Range("H99").Value = 6294.636
$ws.Range("I99").Value = 6294.636
$ws.Range("K99").Value = 6294.636
$ws.Range("M99").Value = -3299.636
$ws.Range("H136").Value = 30304240
$ws.Range("I136").Value = 40000876
$ws.Range("K136").Value = 120002628
$ws.Range("M136").Value = -120000078

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3481.4888
$ws.Range("I134").Value = 1080.7646
$ws.Range("J134").Value = 10901.909
$ws.Range("K134").Value = 3242.2938
$ws.Range("L134").Value = 32705.727
$ws.Range("M134").Value = -707.2937999999999
$ws.Range("N134").Value = -37775.727

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1720.4
$ws.Range("I31").Value = 1771.8572
$ws.Range("J31").Value = 1000
$ws.Range("K31").Value = 1771.8572
$ws.Range("L31").Value = 1000
$ws.Range("M31").Value = -1476.8572
$ws.Range("N31").Value = -1590
$ws.Range("H34").Value = 1720.4
$ws.Range("I34").Value = 1771.8572
$ws.Range("J34").Value = 1000
$ws.Range("K34").Value = 1771.8572
$ws.Range("L34").Value = 1000
$ws.Range("M34").Value = -1569.8572
$ws.Range("N34").Value = -1404
$ws.Range("H58").Value = 1009.34784
$ws.Range("I58").Value = 957.8946999999999
$ws.Range("K58").Value = 957.8946999999999
$ws.Range("M58").Value = -754.8946999999999
$ws.Range("H132").Value = 1557.4445
$ws.Range("I132").Value = 1205.2646
$ws.Range("J132").Value = 2646
$ws.Range("K132").Value = 3615.7938
$ws.Range("L132").Value = 7938
$ws.Range("M132").Value = -1085.7938
$ws.Range("N132").Value = -12998
$ws.Range("H134").Value = 12500998
$ws.Range("I134").Value = 912.2857
$ws.Range("J134").Value = 41667864
$ws.Range("K134").Value = 2736.8571
$ws.Range("L134").Value = 125003592
$ws.Range("M134").Value = -201.8571000000002
$ws.Range("N134").Value = -125008662
$ws.Range("H136").Value = 1009.34784
$ws.Range("I136").Value = 957.8946999999999
$ws.Range("K136").Value = 2873.6841
$ws.Range("M136").Value = -323.6840999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 180000
$ws.Range("J37").Value = 180000
$ws.Range("L37").Value = 540000
$ws.Range("N37").Value = -540224
$ws.Range("H40").Value = 207.72223
$ws.Range("J40").Value = 371.2857
$ws.Range("L40").Value = 1485.1428
$ws.Range("N40").Value = -1623.1428
$ws.Range("H123").Value = 2424.2856
$ws.Range("I123").Value = 1251.25
$ws.Range("J123").Value = 2893.5
$ws.Range("K123").Value = 3753.75
$ws.Range("L123").Value = 8680.5
$ws.Range("M123").Value = -1303.75
$ws.Range("N123").Value = -13580.5
$ws.Range("H129").Value = 16027253
$ws.Range("I129").Value = 66667228
$ws.Range("J129").Value = 3970116
$ws.Range("K129").Value = 200001684
$ws.Range("L129").Value = 11910348
$ws.Range("M129").Value = -199996684
$ws.Range("N129").Value = -11920348
$ws.Range("H130").Value = 2255.0667
$ws.Range("I130").Value = 1030
$ws.Range("J130").Value = 2342.5715
$ws.Range("K130").Value = 3090
$ws.Range("L130").Value = 7027.7145
$ws.Range("M130").Value = 1930
$ws.Range("N130").Value = -17067.7145
$ws.Range("H131").Value = 25003502
$ws.Range("I131").Value = 125000500
$ws.Range("K131").Value = 375001500
$ws.Range("M131").Value = -374996460
$ws.Range("H133").Value = 3552.25
$ws.Range("I133").Value = 2012.5
$ws.Range("J133").Value = 3937.1875
$ws.Range("K133").Value = 6037.5
$ws.Range("L133").Value = 11811.5625
$ws.Range("M133").Value = -977.5
$ws.Range("N133").Value = -21931.5625
$ws.Range("H134").Value = 3986.96
$ws.Range("I134").Value = 1684.4546
$ws.Range("K134").Value = 5053.3638
$ws.Range("M134").Value = 16.63619999999992
$ws.Range("H136").Value = 2832.7144
$ws.Range("I136").Value = 1620
$ws.Range("J136").Value = 4449.6665
$ws.Range("K136").Value = 4860
$ws.Range("L136").Value = 13348.9995
$ws.Range("M136").Value = 240
$ws.Range("N136").Value = -23548.9995
$ws.Range("H137").Value = 28852880
$ws.Range("I137").Value = 57693990
$ws.Range("J137").Value = 11770
$ws.Range("K137").Value = 173081970
$ws.Range("L137").Value = 35310
$ws.Range("M137").Value = -173076870
$ws.Range("N137").Value = -45510
$ws.Range("H138").Value = 1750.45
$ws.Range("I138").Value = 1185.3077
$ws.Range("K138").Value = 3555.9231
$ws.Range("M138").Value = 1584.0769
$ws.Range("H139").Value = 2090.5952
$ws.Range("I139").Value = 2594.0527
$ws.Range("J139").Value = 1674.6957
$ws.Range("K139").Value = 7782.158100000001
$ws.Range("L139").Value = 5024.0871
$ws.Range("M139").Value = -2642.158100000001
$ws.Range("N139").Value = -15304.0871
$ws.Range("H140").Value = 18662.709
$ws.Range("I140").Value = 43692.793
$ws.Range("J140").Value = 2854.2368
$ws.Range("K140").Value = 131078.379
$ws.Range("L140").Value = 8562.7104
$ws.Range("M140").Value = -125898.379
$ws.Range("N140").Value = -18922.7104
$ws.Range("H141").Value = 1926.8462
$ws.Range("I141").Value = 1926.8462
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 5780.5386
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -600.5385999999999
$ws.Range("N141").Value = $null

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 27596.334
$ws.Range("I132").Value = 1266.1052
$ws.Range("J132").Value = 52610.05
$ws.Range("K132").Value = 3798.3156
$ws.Range("L132").Value = 157830.15
$ws.Range("M132").Value = -1268.3156
$ws.Range("N132").Value = -162890.15
$ws.Range("H136").Value = 2238.4614
$ws.Range("I136").Value = 2260
$ws.Range("J136").Value = 2166.6667
$ws.Range("K136").Value = 6780
$ws.Range("L136").Value = 6500.000100000001
$ws.Range("M136").Value = -4230
$ws.Range("N136").Value = -11600.0001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 54011.4
$ws.Range("I28").Value = 30000
$ws.Range("J28").Value = 70019
$ws.Range("K28").Value = 30000
$ws.Range("L28").Value = 70019
$ws.Range("M28").Value = -29652
$ws.Range("N28").Value = -70715
$ws.Range("H132").Value = 1570.42
$ws.Range("I132").Value = 1059.7812
$ws.Range("K132").Value = 3179.3436
$ws.Range("M132").Value = -649.3435999999997
$ws.Range("H136").Value = 723.7143
$ws.Range("I136").Value = 648.7619
$ws.Range("J136").Value = 948.5714
$ws.Range("K136").Value = 1946.2857
$ws.Range("L136").Value = 2845.7142
$ws.Range("M136").Value = 603.7143000000001
$ws.Range("N136").Value = -7945.7142
